$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '65.293.19'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -2.55%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.471.59'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +0.34%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  -0.22%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '554.06'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +1.56%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '179.17'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -3.38%  '
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +4.64%  '
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +0.02%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.634'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -0.35%  '
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +3.50%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '54.04'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -3.13%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0000271'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -0.88%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '9.24'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -2.01%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.015.44'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -0.05%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '18.64'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +2.19%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.466.75'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +0.12%  '
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +0.81%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '11.97'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +1.63%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '65.210.48'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -3.20%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.989'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -1.43%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '415.91'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +2.71%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.05'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +4.28%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '86.09'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +1.82%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '4.27'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +2.24%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '12.92'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +10.23%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '10.84'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -8.87%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.85'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -2.09%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '6.02'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -3.50%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.13'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +5.74%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '30.31'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +0.64%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '6.57'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -4.00%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '608.97'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -9.82%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '11.78'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +1.32%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.110'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -0.16%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '59.10'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +0.15%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.999'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -0.10%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '37.49'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -2.69%  '
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +9.45%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0₃0789'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -4.10%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.367.63'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +10.65%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.380'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -5.53%  '
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -0.05%  '
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -3.36%  '
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -5.71%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.54'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -9.18%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.27'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -0.36%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0413'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -1.08%  '
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -0.99%  '
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +3.07%  '

$ws.Range("B50").NumberFormat = "@"
$ws.Range("B50").Value = 'THORChain'
$ws.Range("C50").NumberFormat = "@"
$ws.Range("C50").Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '8.45'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -3.51%  '

$ws.Range("B51").NumberFormat = "@"
$ws.Range("B51").Value = 'Monero'
$ws.Range("C51").NumberFormat = "@"
$ws.Range("C51").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '137.76'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -1.41%  '

